$d = $word.ActiveDocument

# Anchor: the "Amount balance ... - 1760.0" paragraph (the last populated
# paragraph before the long run of trailing blank paragraphs).
$anchor = $d.Paragraphs.Item(61)
$r = $anchor.Range
$r.Collapse(0)

$dashes = "---------------------------------------------------------------"

$text = "" `
  + "`r" `
  + "`r" `
  + "SUN OCT 1 10:52:54 PDT 2017" `
  + "`r" `
  + "Person Name`t`t`t`t- YASHODHA" `
  + "`r" `
  + $dashes `
  + "`r" `
  + "Item Name`t`t`t`t- CARROT" `
  + "`r" `
  + "Amount Received`t`t`t- 1760" `
  + "`r" `
  + "Amount Received mode`t`t- CASH AND CLEARD" `
  + "`r" `
  + "`r"

$r.InsertAfter($text)

# After the insertion, the new content occupies paragraphs 62..71 (10
# paragraphs), and the document's original trailing blank paragraphs
# resume at 72.
#   62: (blank, bold)
#   63: (blank, bold)
#   64: SUN OCT 1 10:52:54 PDT 2017
#   65: Person Name ... - YASHODHA
#   66: -------------------------- (63 dashes)
#   67: Item Name ... - CARROT
#   68: Amount Received ... - 1760   (red)
#   69: Amount Received mode ... - CASH AND CLEARD
#   70: (blank)
#   71: (blank, bold)

for ($i = 62; $i -le 71; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.Name = "Courier New"
}

$d.Paragraphs.Item(62).Range.Font.Bold = $true
$d.Paragraphs.Item(63).Range.Font.Bold = $true
$d.Paragraphs.Item(68).Range.Font.Color = 255
$d.Paragraphs.Item(71).Range.Font.Bold = $true

Write-Host ("Paragraphs.Count=" + $d.Paragraphs.Count)
for ($i = 58; $i -le 74; $i++) {
    Write-Host ($i.ToString() + ": [" + $d.Paragraphs.Item($i).Range.Text + "]")
}
